$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff": the 6 files that were "Ready for handoff"
# (rows 4-9 on every sheet -- row 10 already has downstream handback data and
# is left untouched) get a fresh handoff pass: Priority is stamped "ht" and
# the Latest Handoff Datetime / Latest HO Xliff Generate Date are refreshed.

# --- Overview sheet: column G = "Latest HO Xliff Generate Date" ---
$wsOverview = $wb.Sheets.Item("Overview")
for ($r = 4; $r -le 9; $r++) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-10-21 01:37:12"
}

# --- zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime" ---
$wsZhCn = $wb.Sheets.Item("zh-cn")
for ($r = 4; $r -le 9; $r++) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsZhCn.Cells.Item($r, 8).Value = "2016-10-21 01:37:00"
}

# --- de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime" ---
$wsDeDe = $wb.Sheets.Item("de-de")
for ($r = 4; $r -le 9; $r++) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 8).Value = "2016-10-21 01:37:12"
}
